$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update data values
$ws.Range("B2").Value = 58.82
$ws.Range("D2").Value = 34

$ws.Range("B6").Value = 56.25

$ws.Range("C7").Value = 63.79
$ws.Range("E7").Value = 116

$ws.Range("C8").Value = 72.34
$ws.Range("E8").Value = 94

$ws.Range("C9").Value = 71.29
$ws.Range("E9").Value = 101

# Update the active selection to E10 (matches activeCell/sqref change in diff)
$ws.Range("E10").Select()
